$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45805.47588245371
$ws.Range("B2").Value = 0.03337472496505753
$ws.Range("C2").Value = "Ổn định giá, Tăng trưởng dài hạn, tieu chi khac"
$ws.Range("D2").Value = "Vàng, Ngoại tệ, phuong an khac"
$ws.Range("E2").Value = "Vàng (0.3333), Ngoại tệ (0.3333), phuong an khac (0.3333)"

# Row 3
$ws.Range("A3").Value = 45798.58740990741
$ws.Range("B3").Value = 0.03337472496505753
$ws.Range("C3").Value = "Ổn định giá, Tính thanh khoản, Tăng trưởng dài hạn"
$ws.Range("D3").Value = "Ngoại tệ, Cổ phiếu, Vàng"
$ws.Range("E3").Value = "Ngoại tệ (0.3333), Cổ phiếu (0.3333), Vàng (0.3333)"

# Row 4
$ws.Range("A4").Value = 45798.58724861111
$ws.Range("B4").Value = 0.03337472496505753
$ws.Range("C4").Value = "Ổn định giá, Tính thanh khoản, Tăng trưởng dài hạn"
$ws.Range("D4").Value = "Ngoại tệ, Cổ phiếu, Vàng"
$ws.Range("E4").Value = "Ngoại tệ (0.3333), Cổ phiếu (0.3333), Vàng (0.3333)"

# Row 5
$ws.Range("A5").Value = 45798.57278447917
$ws.Range("B5").Value = 0.06390856373847034
$ws.Range("C5").Value = "Ổn định giá, Tăng trưởng dài hạn, Tính thanh khoản"
$ws.Range("D5").Value = "Ngoại tệ, Vàng, Cổ phiếu"
$ws.Range("E5").Value = "Ngoại tệ (0.3333), Vàng (0.2884), Cổ phiếu (0.2196)"

# Row 6
$ws.Range("A6").Value = 45798.56284561343
$ws.Range("B6").Value = 0.03337472496505753
$ws.Range("C6").Value = "Ổn định giá, Tăng trưởng dài hạn, Tính thanh khoản"
$ws.Range("D6").Value = "Ngoại tệ, Vàng, Cổ phiếu"
$ws.Range("E6").Value = "Cổ phiếu (0.3832), Ngoại tệ (0.3333), Vàng (0.2879)"

# Row 7 (new)
$ws.Range("A7").Value = 45798.55942474537
$ws.Range("A7").NumberFormat = $ws.Range("A6").NumberFormat
$ws.Range("B7").Value = 0.03337472496505753
$ws.Range("C7").Value = "Ổn định giá, Tăng trưởng dài hạn, Khả năng chống lạm phát"
$ws.Range("D7").Value = "Vàng, Ngoại tệ, Cổ phiếu"
$ws.Range("E7").Value = "Vàng (0.3333), Ngoại tệ (0.2771), Cổ phiếu (0.2406)"
